$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the weekly data between row 3 and row 4 (Fecha, Volumen, Precio minimo,
# Precio promedio ponderado, Precio $/Kg), matching the target diff.

$ws.Range("D3").Value = 44804
$ws.Range("J3").Value = 50
$ws.Range("K3").Value = 9500
$ws.Range("M3").Value = 9750
$ws.Range("P3").Value = 542

$ws.Range("D4").Value = 44714
$ws.Range("J4").Value = 80
$ws.Range("K4").Value = 9000
$ws.Range("M4").Value = 9500
$ws.Range("P4").Value = 528
